$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "Inflammatory-Mac"
$ws.Range("B2").Value = "Il12b"
$ws.Range("C2").Value = "Il12rb2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 1.284858666666667
$ws.Range("H2").Value = 3.854576
$ws.Range("I2").Value = 0.6825120125588942
$ws.Range("J2").Value = 0.6825120125588942
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.1320293333333333
$ws.Range("N2").Value = 0.396088
$ws.Range("O2").Value = 0.02700478969442551
$ws.Range("P2").Value = 0.02700478969442551
$ws.Range("Q2").Value = 0.1696390331875556
$ws.Range("R2").Value = 1.526751298688
$ws.Range("S2").Value = 0.01843109336307204
$ws.Range("T2").Value = 0.01843109336307204

# Row 3
$ws.Range("A3").Value = "Inflammatory-Mac"
$ws.Range("B3").Value = "Il12b"
$ws.Range("C3").Value = "Il12rb2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 1.284858666666667
$ws.Range("H3").Value = 3.854576
$ws.Range("I3").Value = 0.6825120125588942
$ws.Range("J3").Value = 0.6825120125588942
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.5252536666666666
$ws.Range("N3").Value = 1.575761
$ws.Range("O3").Value = 0.1074334350287755
$ws.Range("P3").Value = 0.1074334350287755
$ws.Range("Q3").Value = 0.6748767258151112
$ws.Range("R3").Value = 6.073890532336001
$ws.Range("S3").Value = 0.07332460995760477
$ws.Range("T3").Value = 0.07332460995760477

# Row 4
$ws.Range("A4").Value = "Inflammatory-Mac"
$ws.Range("B4").Value = "Il12b"
$ws.Range("C4").Value = "Il12rb2"
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 1.284858666666667
$ws.Range("H4").Value = 3.854576
$ws.Range("I4").Value = 0.6825120125588942
$ws.Range("J4").Value = 0.6825120125588942
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.687203666666667
$ws.Range("N4").Value = 5.061611
$ws.Range("O4").Value = 0.3450943744066743
$ws.Range("P4").Value = 0.3450943744066743
$ws.Range("Q4").Value = 2.167818253548445
$ws.Range("R4").Value = 19.510364281936
$ws.Range("S4").Value = 0.2355310559990518
$ws.Range("T4").Value = 0.2355310559990518

# Row 5
$ws.Range("A5").Value = "Inflammatory-Mac"
$ws.Range("B5").Value = "Il12b"
$ws.Range("C5").Value = "Il12rb2"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 1.284858666666667
$ws.Range("H5").Value = 3.854576
$ws.Range("I5").Value = 0.6825120125588942
$ws.Range("J5").Value = 0.6825120125588942
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.505857
$ws.Range("N5").Value = 1.517571
$ws.Range("O5").Value = 0.1034661128369428
$ws.Range("P5").Value = 0.1034661128369428
$ws.Range("Q5").Value = 0.6499547505440001
$ws.Range("R5").Value = 5.849592754896
$ws.Range("S5").Value = 0.0706168649039875
$ws.Range("T5").Value = 0.0706168649039875

# Row 6
$ws.Range("A6").Value = "Inflammatory-Mac"
$ws.Range("B6").Value = "Il12b"
$ws.Range("C6").Value = "Il12rb2"
$ws.Range("D6").Value = "Resolving-Mac"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 1.284858666666667
$ws.Range("H6").Value = 3.854576
$ws.Range("I6").Value = 0.6825120125588942
$ws.Range("J6").Value = 0.6825120125588942
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 2.038764333333333
$ws.Range("N6").Value = 6.116293000000001
$ws.Range("O6").Value = 0.4170012880331818
$ws.Range("P6").Value = 0.4170012880331818
$ws.Range("Q6").Value = 2.619524022974223
$ws.Range("R6").Value = 23.575716206768
$ws.Range("S6").Value = 0.284608388335178
$ws.Range("T6").Value = 0.284608388335178

# Row 7
$ws.Range("A7").Value = "Resolving-Mac"
$ws.Range("B7").Value = "Il12b"
$ws.Range("C7").Value = "Il12rb2"
$ws.Range("D7").Value = "ECs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.597685
$ws.Range("H7").Value = 1.793055
$ws.Range("I7").Value = 0.3174879874411058
$ws.Range("J7").Value = 0.3174879874411058
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.1320293333333333
$ws.Range("N7").Value = 0.396088
$ws.Range("O7").Value = 0.02700478969442551
$ws.Range("P7").Value = 0.02700478969442551
$ws.Range("Q7").Value = 0.07891195209333333
$ws.Range("R7").Value = 0.71020756884
$ws.Range("S7").Value = 0.008573696331353471
$ws.Range("T7").Value = 0.008573696331353471

# Row 8
$ws.Range("A8").Value = "Resolving-Mac"
$ws.Range("B8").Value = "Il12b"
$ws.Range("C8").Value = "Il12rb2"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.597685
$ws.Range("H8").Value = 1.793055
$ws.Range("I8").Value = 0.3174879874411058
$ws.Range("J8").Value = 0.3174879874411058
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.5252536666666666
$ws.Range("N8").Value = 1.575761
$ws.Range("O8").Value = 0.1074334350287755
$ws.Range("P8").Value = 0.1074334350287755
$ws.Range("Q8").Value = 0.3139362377616666
$ws.Range("R8").Value = 2.825426139855
$ws.Range("S8").Value = 0.03410882507117075
$ws.Range("T8").Value = 0.03410882507117075

# Row 9
$ws.Range("A9").Value = "Resolving-Mac"
$ws.Range("B9").Value = "Il12b"
$ws.Range("C9").Value = "Il12rb2"
$ws.Range("D9").Value = "Inflammatory-Mac"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.597685
$ws.Range("H9").Value = 1.793055
$ws.Range("I9").Value = 0.3174879874411058
$ws.Range("J9").Value = 0.3174879874411058
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1.687203666666667
$ws.Range("N9").Value = 5.061611
$ws.Range("O9").Value = 0.3450943744066743
$ws.Range("P9").Value = 0.3450943744066743
$ws.Range("Q9").Value = 1.008416323511667
$ws.Range("R9").Value = 9.075746911605
$ws.Range("S9").Value = 0.1095633184076225
$ws.Range("T9").Value = 0.1095633184076225

# Row 10
$ws.Range("A10").Value = "Resolving-Mac"
$ws.Range("B10").Value = "Il12b"
$ws.Range("C10").Value = "Il12rb2"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.597685
$ws.Range("H10").Value = 1.793055
$ws.Range("I10").Value = 0.3174879874411058
$ws.Range("J10").Value = 0.3174879874411058
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.505857
$ws.Range("N10").Value = 1.517571
$ws.Range("O10").Value = 0.1034661128369428
$ws.Range("P10").Value = 0.1034661128369428
$ws.Range("Q10").Value = 0.302343141045
$ws.Range("R10").Value = 2.721088269405
$ws.Range("S10").Value = 0.03284924793295535
$ws.Range("T10").Value = 0.03284924793295535

# Row 11
$ws.Range("A11").Value = "Resolving-Mac"
$ws.Range("B11").Value = "Il12b"
$ws.Range("C11").Value = "Il12rb2"
$ws.Range("D11").Value = "Resolving-Mac"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.597685
$ws.Range("H11").Value = 1.793055
$ws.Range("I11").Value = 0.3174879874411058
$ws.Range("J11").Value = 0.3174879874411058
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 2.038764333333333
$ws.Range("N11").Value = 6.116293000000001
$ws.Range("O11").Value = 0.4170012880331818
$ws.Range("P11").Value = 0.4170012880331818
$ws.Range("Q11").Value = 1.218538860568333
$ws.Range("R11").Value = 10.966849745115
$ws.Range("S11").Value = 0.1323928996980038
$ws.Range("T11").Value = 0.1323928996980038
